$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the stray _GoBack bookmark that currently sits after
#    "and literary " (it will be re-created later at the location of
#    the real last edit, i.e. after the corrected Verhagen affiliation).
# ------------------------------------------------------------------
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# ------------------------------------------------------------------
# 2) Split the "Philip Verhagen, Vryje University" run right after the
#    leading "P" so the correction below lands in its own run (mirrors
#    how Word splits a run when you start retyping partway through it).
#    We do this with a throwaway bookmark purely to force the split;
#    it gets removed again immediately after.
# ------------------------------------------------------------------
$rngP = $d.Content
$rngP.Find.Execute("Philip Verhagen, Vryje University", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $rngP.Duplicate
$splitPoint.SetRange($rngP.Start + 1, $rngP.Start + 1)
$d.Bookmarks.Add("TEMP_SPLIT", $splitPoint)

# ------------------------------------------------------------------
# 3) Fix the affiliation text: "hilip Verhagen, Vryje University"
#    becomes "hilip Verhagen, Vrije Universiteit Amsterdam".
# ------------------------------------------------------------------
$tail = $d.Content
$tail.Find.Execute("hilip Verhagen, Vryje University", $true, $false, $false, $false, $false, $true, 1, $false, "hilip Verhagen, Vrije Universiteit Amsterdam", 2)

$tb = $d.Bookmarks("TEMP_SPLIT")
$tb.Delete()

# ------------------------------------------------------------------
# 4) Re-create _GoBack right after "Amsterdam" (the end of the edit),
#    matching Word's habit of marking the spot of the last change.
#    A bookmark collapsed exactly at a paragraph's final position gets
#    snapped to span the whole paragraph, so we anchor it just before a
#    throwaway character appended after "Amsterdam" and then remove
#    that character, leaving the bookmark collapsed in the right spot.
# ------------------------------------------------------------------
$rngFind = $d.Content
$rngFind.Find.Execute("Amsterdam", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPt = $rngFind.Duplicate
$endPt.Collapse(0)
$zStart = $endPt.Start
$endPt.InsertAfter("Z")

$beforeZ = $d.Range($zStart, $zStart)
$d.Bookmarks.Add("_GoBack", $beforeZ)

$zRange = $d.Range($zStart, $zStart + 1)
$zRange.Delete()
